# Hortaliza, Vega Modelo de Temuco - Brocoli: add a new weekly price record
# (3 new rows) ahead of the existing data, shifting the rest down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 623. This pushes the current rows
# 623-719 down to 626-722, preserving all of their data/formatting.
$ws.Range("A623:A625").EntireRow.Insert()

# Columns that are constant for every data row in this sheet.
$constA = 10
$constB = "Vega Modelo de Temuco"
$constC = "La Araucanía"
$constE = 9
$constF = 100112023
$constG = "Brócoli"
$constH = "Sin especificar"
$constN = "`$/unidad"
$constQ = 1
$constR = "Hortaliza"

function Set-DataRow {
    param($row, $D, $I, $J, $K, $L, $M, $O, $P)
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $K
    $ws.Cells.Item($row, 12).Value = $L
    $ws.Cells.Item($row, 13).Value = $M
    $ws.Cells.Item($row, 14).Value = $constN
    $ws.Cells.Item($row, 15).Value = $O
    $ws.Cells.Item($row, 16).Value = $P
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
}

# New records (dated 2023-07-20 / serial 45127)
Set-DataRow 623 45127 "Primera" 2200 950  1000 973  "Región Metropolitana" 973
Set-DataRow 624 45127 "Primera" 100  1200 1200 1200 "Región de O'Higgins"  1200
Set-DataRow 625 45127 "Primera" 4000 1000 1000 1000 "Región del Maule"     1000
